# Insert two new weekly price records for "Cilantro" at Terminal Hortofrutícola
# Agro Chillán, right before the existing row 57 (2022-11-03 entry), shifting
# all subsequent rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("57:58").Insert()

# New row 57: Primera quality, 2022-11-18 (serial 44883)
$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 44883
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = 100112040
$ws.Cells.Item(57, 7).Value = "Cilantro"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 300
$ws.Cells.Item(57, 11).Value = 600
$ws.Cells.Item(57, 12).Value = 700
$ws.Cells.Item(57, 13).Value = 650
$ws.Cells.Item(57, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(57, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(57, 16).Value = 650
$ws.Cells.Item(57, 17).Value = 1
$ws.Cells.Item(57, 18).Value = "Hortaliza"

# New row 58: Segunda quality, 2022-11-18 (serial 44883)
$ws.Cells.Item(58, 1).Value = 7
$ws.Cells.Item(58, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value = "Ñuble"
$ws.Cells.Item(58, 4).Value = 44883
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 100112040
$ws.Cells.Item(58, 7).Value = "Cilantro"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Segunda"
$ws.Cells.Item(58, 10).Value = 200
$ws.Cells.Item(58, 11).Value = 500
$ws.Cells.Item(58, 12).Value = 500
$ws.Cells.Item(58, 13).Value = 500
$ws.Cells.Item(58, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(58, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(58, 16).Value = 500
$ws.Cells.Item(58, 17).Value = 1
$ws.Cells.Item(58, 18).Value = "Hortaliza"
